# "Final tweaks before test"
#
# Substantive content changes captured by the xml diff:
#   1. H5's HYPERLINK formula had a slash dropped from the URL
#      (https://shorturl.at/ALPX4 -> https:/shorturl.at/ALPX4).
#   2. The active selection moved from E14 to H5 (the cell that was
#      just edited).
#
# (Everything else in the diff -- fileVersion/build numbers, GUIDs,
# the absPath/revisionPtr bookkeeping, theme display-name relabelling,
# the numFmt id renumbering, the new empty xl/persons/person.xml part,
# x14ac:dyDescent/knownFonts noise, and the window geometry -- is the
# kind of silent re-serialization churn Excel performs on every save
# and isn't driven by any explicit user action, so it isn't reproduced
# here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd hyperlink in H5 (dropped slash after "https:").
$ws.Range("H5").Formula = '=HYPERLINK("https:/shorturl.at/ALPX4")'

# Leave the selection on H5, matching the saved view state.
$ws.Range("H5").Select()
